# "distribution add information div"
#
# Row 767 (the last data row) had #N/A errors in the Latitude / Longitude
# columns (E/F) left over from a lookup that didn't resolve. Clear those
# error values out (keeping the existing cell style) and bring the view
# down to the bottom of the sheet, selecting the now-empty trailing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the stray #N/A values out of E767 / F767 without touching their
# number format / style.
$ws.Range("E767").ClearContents()
$ws.Range("F767").ClearContents()

# Scroll the frozen-pane view down near the bottom of the data and select
# the last row (A768:XFD768 - a whole-row selection one past the last
# populated row, as Excel leaves it after selecting a row header).
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 760

$ws.Rows.Item(768).Select()
